$d = $word.ActiveDocument

# The paragraph currently holds 15 filler "words" (Aaaa, bbbb, cccc, ...)
# each wrapped in <w:proofErr> spell-check markers, joined by single
# spaces, ending in "=". Replace the whole paragraph's content with the
# real sentence, split across seven runs (all keeping the original
# Consolas run font), add a matching paragraph-mark run font in pPr/rPr,
# and drop the now-irrelevant proofErr spell-check bookmarks.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t>Map shows a</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t xml:space="preserve"> 12-</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t xml:space="preserve">route plan </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t xml:space="preserve">for </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t>seeing</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t xml:space="preserve"> nearly all listed areas</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/></w:rPr><w:t>.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@

$d.Content.InsertXML($xml)
